$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RegisteredNumbers (D), Ideal Number (E), Difference (F)
$ws.Range("D2").Value = 68
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 58

$ws.Range("D3").Value = 52
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 42

$ws.Range("D4").Value = 65
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 55

$ws.Range("D5").Value = 55
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 45

$ws.Range("D6").Value = 34
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 29

$ws.Range("D7").Value = 26
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 21

# Remove the "Re-balanced Number" column values for rows 2-7
$ws.Range("G2:G7").ClearContents()
